# Add the new "compliant_wall" worksheet after the last existing sheet
# ("Turbulent convection NREL clust"), matching sheetId=3 / rId3 order.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "compliant_wall"

# Header row (row 1) - mirrors the header layout used on the
# "Turbulent convection NREL clust" sheet, but with the compliant-wall
# specific coefficients Cm / Cd / Ck in place of Pr / A (dP/dx) / omega (dP/dx).
$new.Range("A1").Value = "NREL job ID"
$new.Range("B1").Value = "Re_tau"
$new.Range("C1").Value = "Cm"
$new.Range("D1").Value = "Cd"
$new.Range("E1").Value = "Ck"
$new.Range("F1").Value = "Lx"
$new.Range("G1").Value = "Ly"
$new.Range("H1").Value = "Lz"
$new.Range("I1").Value = "Nx"
$new.Range("J1").Value = "Ny"
$new.Range("K1").Value = "Nz"
$new.Range("L1").Value = "initial_dt"
$new.Range("M1").Value = "ntask"
$new.Range("N1").Value = "partition"
$new.Range("O1").Value = "mem"
$new.Range("P1").Value = "time"
$new.Range("Q1").Value = "Note"

# Data row (row 2)
$new.Range("B2").Value = 180
$new.Range("C2").Value = 2
$new.Range("D2").Value = -2.93
$new.Range("E2").Value = 28859
$new.Range("F2").Value = "4pi"
$new.Range("G2").Value = 2
$new.Range("H2").Value = "2pi"
$new.Range("I2").Value = 192
$new.Range("J2").Value = 258
$new.Range("K2").Value = 160
$new.Range("L2").Value = 0.0000001
$new.Range("L2").NumberFormat = "0.00E+00"
$new.Range("M2").Value = 96
$new.Range("N2").Value = "standard"

# Restore the previous sheet's lingering selection (whole first two rows,
# scrolled down toward the bottom of its data) before the new sheet took
# over as the active tab.
$ws2 = $wb.Worksheets.Item("Turbulent convection NREL clust")
$ws2.Range("A1:XFD2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49

# Re-activate the new sheet so it ends up as the selected/visible tab, with
# the cursor left on E2 (the last cell filled in on the new row).
$new.Activate() | Out-Null
$new.Range("E2").Select() | Out-Null
